$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (company 2) ---
$ws.Range("D2").Value = 0.1845
$ws.Range("E2").Value = 0.148
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 37.84
$ws.Range("L2").Value = 0.3158597662771285
$ws.Range("M2").Value = 6.35
$ws.Range("N2").Value = 0.007773289264291835
$ws.Range("O2").Value = 0.1678118393234672
$ws.Range("P2").Value = 6.35
$ws.Range("Q2").Value = 0.007773289264291835
$ws.Range("R2").Value = 0.1678118393234672
$ws.Range("U2").Value = 2212.7
$ws.Range("V2").Value = 2.708654670094258
$ws.Range("W2").Value = 0.1317822641509434
$ws.Range("X2").Value = 0.04850390130502357
$ws.Range("Y2").Value = 0.08327836284591983
$ws.Range("Z2").Value = -0.09436711802191397
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.03727689281575863
$ws.Range("AC2").Value = -0.03727689281575863
$ws.Range("AD2").Value = 738.6
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 738.6
$ws.Range("AG2").Value = -1474.1
$ws.Range("AH2").Value = 0.474831243972999
$ws.Range("AI2").Value = 0.663611859838275
$ws.Range("AJ2").Value = 2.243000608642727
$ws.Range("AK2").Value = 1.340456488133127
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# --- Row 3 (Coop Pank AS) ---
$ws.Range("B3").Value = "Coop Pank AS (TLSE:CPA1T)"
$ws.Range("D3").Value = 0.194
$ws.Range("E3").Value = 0.178
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 7.54
$ws.Range("L3").Value = 0.2654929577464789
$ws.Range("U3").Value = 191.5
$ws.Range("V3").Value = 1.478764478764479
$ws.Range("W3").Value = 0.12064
$ws.Range("X3").Value = 0.04057935693798563
$ws.Range("Y3").Value = 0.08006064306201437
$ws.Range("Z3").Value = -1.663737551259519
$ws.Range("AA3").Value = -0
$ws.Range("AB3").Value = 0.03639896689941027
$ws.Range("AC3").Value = -0.03639896689941027
$ws.Range("AD3").Value = 33.7
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 33.7
$ws.Range("AG3").Value = -157.8
$ws.Range("AH3").Value = 0.2064950980392157
$ws.Range("AI3").Value = 0.2333795013850416
$ws.Range("AJ3").Value = 5.575971731448761
$ws.Range("AK3").Value = 3.350318471337579
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()

# --- Row 4 (AS LHV Group) ---
$ws.Range("D4").Value = 0.175
$ws.Range("E4").Value = 0.118
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 30.3
$ws.Range("L4").Value = 0.3315098468271335
$ws.Range("M4").Value = 6.35
$ws.Range("N4").Value = 0.009237707302880418
$ws.Range("O4").Value = 0.2095709570957096
$ws.Range("P4").Value = 6.35
$ws.Range("Q4").Value = 0.009237707302880418
$ws.Range("R4").Value = 0.2095709570957096
$ws.Range("U4").Value = 2021.2
$ws.Range("V4").Value = 2.94035496072156
$ws.Range("W4").Value = 0.1429245283018868
$ws.Range("X4").Value = 0.05642844567206151
$ws.Range("Y4").Value = 0.08649608262982528
$ws.Range("Z4").Value = -0.0729775478266424
$ws.Range("AA4").Value = -0
$ws.Range("AB4").Value = 0.03815481873210701
$ws.Range("AC4").Value = -0.03815481873210701
$ws.Range("AD4").Value = 704.9
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 704.9
$ws.Range("AG4").Value = -1316.3
$ws.Range("AH4").Value = 0.5062845651080945
$ws.Range("AI4").Value = 0.7277513937641958
$ws.Range("AJ4").Value = 2.093019557958339
$ws.Range("AK4").Value = 1.25052251567547
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
